# Insert a new data row at row 396 (pushing existing rows 396-427 down to 397-428)
# for the "Fruta, Vega Central Mapocho de Santiago - Mango" weekly price sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before current row 396; all rows from 396 onward shift down by one.
$ws.Rows.Item(396).Insert()

# Populate the newly inserted row 396 with the new weekly record.
$ws.Range("A396").Value = 9
$ws.Range("B396").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C396").Value = "Metropolitana"
$ws.Range("D396").Value = 44714
$ws.Range("E396").Value = 13
$ws.Range("F396").Value = "Fruta"
$ws.Range("G396").Value = 100108
$ws.Range("H396").Value = "Tropicales y subtropicales"
$ws.Range("I396").Value = 100108002
$ws.Range("J396").Value = "Mango"
$ws.Range("K396").Value = "Sin especificar"
$ws.Range("L396").Value = "Primera"
$ws.Range("M396").Value = 470
$ws.Range("N396").Value = 9000
$ws.Range("O396").Value = 10000
$ws.Range("P396").Value = 9532
$ws.Range("Q396").Value = "`$/bandeja 4 kilos"
$ws.Range("R396").Value = "Brasil"
$ws.Range("S396").Value = 2383
$ws.Range("T396").Value = 4
